$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "Sembol",
    "GOOGL",
    "GOOG",
    "MSFT",
    "AMZN",
    "META",
    "JPM",
    "LLY",
    "V",
    "ORCL",
    "BAC",
    "MA",
    "BABA",
    "AMD",
    "PLTR",
    "ABBV",
    "NFLX",
    "JPM-PD",
    "JPM-PC",
    "UNH",
    "MS",
    "BAC-PK",
    "BAC-PL",
    "C",
    "NVO",
    "SAP",
    "IBM",
    "WFC",
    "BAC-PB",
    "AXP",
    "TMO",
    "CRM",
    "WFC-PY",
    "DIS",
    "BA",
    "WFC-PL",
    "ISRG",
    "SAN",
    "SCHW",
    "BX",
    "APP",
    "SHOP",
    "ACN",
    "BLK",
    "UBER",
    "DHR",
    "TJX",
    "BKNG",
    "QCOM",
    "HDB",
    "SPGI",
    "INTU",
    "UBS",
    "PDD",
    "BBVA",
    "NOW",
    "COF",
    "BSX",
    "SONY",
    "PANW",
    "ADBE",
    "WFC-PC",
    "VRTX",
    "CRWD",
    "MELI",
    "KKR",
    "CVS",
    "ADP",
    "CEG",
    "CVNA",
    "ICE",
    "GD",
    "SNPS",
    "HOOD",
    "NKE",
    "MCO",
    "BCS",
    "MRSH",
    "DASH",
    "NU",
    "NTES",
    "ELV",
    "ORLY",
    "MS-PK",
    "MS-PI",
    "ABNB",
    "TDG",
    "DB",
    "MS-PF",
    "MS-PE",
    "APO",
    "GM",
    "INFY",
    "SE",
    "USB-PP",
    "AON",
    "SNOW",
    "RELX",
    "NWG",
    "LHX",
    "AJG",
    "DUK-PA",
    "RKT",
    "CTA-PB",
    "ADSK",
    "COIN",
    "NDAQ",
    "IDXX",
    "TRI",
    "BIDU",
    "PYPL",
    "CMG",
    "RBLX",
    "EA",
    "SCHW-PD",
    "WDAY",
    "CBRE",
    "EW",
    "ARES",
    "AXON",
    "ALNY",
    "ROK",
    "AMP",
    "CTA-PA",
    "HEI",
    "MSTR",
    "MSCI",
    "TTWO",
    "SPG-PJ",
    "ROP",
    "JD",
    "EBAY",
    "RKLB",
    "MET-PA",
    "EL",
    "CTSH",
    "TCOM",
    "LVS",
    "IQV",
    "PUK",
    "RDDT",
    "CPRT",
    "XYZ",
    "HEI-A",
    "ALC",
    "MLM",
    "ASTS",
    "A",
    "PRU",
    "PAYX",
    "MDLN",
    "RMD",
    "FICO",
    "VEEV",
    "GEHC",
    "FISV",
    "TEAM",
    "RYAAY",
    "CPNG",
    "CLS",
    "SYM",
    "APO-PA",
    "XYL",
    "SATS",
    "ZS",
    "INSM",
    "NTRA",
    "RJF",
    "MDB",
    "EXPE",
    "ESLT",
    "PSA-PK",
    "ALL-PH",
    "SOFI",
    "ALL-PB",
    "FOXA",
    "HUM",
    "WTW",
    "FIS",
    "FOX",
    "VRSK",
    "FLUT",
    "MTD",
    "SYF",
    "DXCM",
    "LPLA",
    "NTRS",
    "CBOE",
    "STLA",
    "CSGP",
    "BAP",
    "ALAB",
    "WIT",
    "HBAN",
    "BRO",
    "EXE",
    "EFX",
    "FSLR",
    "STE",
    "AWK",
    "OMC",
    "VLTO",
    "DLR-PK",
    "CINF",
    "FCNCA",
    "CW",
    "BR",
    "LDOS",
    "SQM",
    "AXIA-PC",
    "ILMN",
    "VRSN",
    "TPG",
    "TROW",
    "WAT",
    "NBIS",
    "LULU",
    "OWL",
    "CNC",
    "FUTU",
    "AFRM",
    "DLR-PJ",
    "CYBR",
    "FWONK",
    "PSLV",
    "FWONA",
    "CG",
    "RBA",
    "UTHR",
    "GPN",
    "SSNC",
    "GMAB",
    "PFG",
    "TWLO",
    "Q",
    "INCY",
    "HL",
    "CHKP",
    "PTC",
    "LTM",
    "TOST",
    "GIB",
    "RIVN",
    "PODD",
    "TYL",
    "RVMD",
    "BWXT",
    "KTOS",
    "MRNA",
    "HIG-PG",
    "GRAB",
    "IOT",
    "U",
    "DKS",
    "HPQ",
    "CRCL",
    "FITBI",
    "XPEV",
    "IT",
    "PSNYW",
    "ALLY",
    "PNR",
    "PINS",
    "WST",
    "HUBS",
    "NWS",
    "IREN",
    "ZG",
    "FN",
    "TRMB",
    "MEDP",
    "JLL",
    "Z",
    "TRU",
    "TTD",
    "NLY",
    "HII",
    "TLN",
    "KSPI",
    "ROKU",
    "GEN",
    "NWSA",
    "IONQ",
    "DKNG",
    "AVAV",
    "KEY-PI",
    "EMA",
    "BBIO",
    "GH",
    "MLI",
    "HMY",
    "PFGC",
    "ULS",
    "GDDY",
    "ERIE",
    "ARCC",
    "ICLR",
    "FIG",
    "RGC",
    "PNFP",
    "KRMN",
    "W",
    "CACI",
    "PEN",
    "EVR",
    "CELH",
    "DPZ",
    "BBY",
    "EMBJ",
    "EQH",
    "GWRE",
    "RBRK",
    "HRL",
    "FIGR",
    "NLY-PG",
    "NLY-PF",
    "NBIX",
    "RVTY",
    "JKHY",
    "RYAN",
    "SF",
    "PSKY",
    "CHWY",
    "UNM",
    "SNAP",
    "JEF",
    "OKLO",
    "HLI",
    "IVZ",
    "AGNC",
    "BMNR",
    "GLXY",
    "GMED",
    "DT",
    "ACGLO",
    "TXRH",
    "DOC",
    "SMMT",
    "RMBS",
    "EPAM",
    "JOBY",
    "CMA",
    "TEM",
    "FHN",
    "AIZ",
    "EXEL",
    "NTNX",
    "BSY",
    "DOCU",
    "BXP",
    "STN",
    "MDGL",
    "BAH",
    "MICC",
    "MP",
    "QGEN",
    "WTRG",
    "SARO",
    "CRL",
    "UHAL",
    "MOH",
    "DRS",
    "MANH",
    "FDS",
    "AFG",
    "CART",
    "SEIC",
    "TECH",
    "CAE",
    "YMM",
    "PCOR",
    "KLAR",
    "JAZZ",
    "CHYM",
    "UHAL-B",
    "SAIL",
    "SANM",
    "GAP",
    "COMP",
    "BIO-B",
    "TTMI",
    "SOLS",
    "ARE",
    "BROS",
    "REXR",
    "AAL",
    "AYI",
    "UWMC",
    "GTLS",
    "ARWR",
    "RGEN",
    "DOX",
    "STEP",
    "MORN",
    "AMG",
    "UEC",
    "LUMN",
    "GGAL",
    "QBTS",
    "TTAN",
    "RZB",
    "AGNCM",
    "AGNCN",
    "PEGA",
    "VIPS",
    "UGI",
    "AMTM",
    "PL"
)

$n = $values.Count
for ($i = 0; $i -lt $n; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
